$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column B header: "year" -> "season_ending_year"
$ws.Range("B1").Value = "season_ending_year"

# Add birth_year (1935) values for all data rows (rows 2-14) in column E
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 5).Value = 1935
}

# Add new column AM: "calendar_year"
$ws.Range("AL1").Copy()
$ws.Range("AM1").PasteSpecial(-4122)
$ws.Range("AM1").Value = "calendar_year"

$calendarYears = @{
    2  = 1969
    3  = 1968
    4  = 1967
    5  = 1966
    6  = 1965
    7  = 1964
    8  = 1963
    9  = 1962
    10 = 1961
    11 = 1960
    12 = 1959
    13 = 1958
    14 = 1957
}

foreach ($r in $calendarYears.Keys) {
    $ws.Cells.Item($r, 39).Value = $calendarYears[$r]
}
